# Update workbook for release "mines - version 1.0.0 (Feb 3 2026)"
# 1) Refresh the version strings on the "About" sheet.
# 2) Remove three retired point features (ids M1439.P2, M1439.P5, M1439.P7)
#    from the "Boundaries and methane sources" sheet, and refresh the
#    build_version column for the rows that remain.

$wb = $excel.ActiveWorkbook

$oldVersion = "mines - January 30 (built on February 02 2026 12.49.33 EST)"
$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- Sheet: About ---------------------------------------------------------
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newVersion

$about.Range("A6").Value = "Recommended Citation:  ""Global Energy Monitor, Coal mine boundaries and methane sources for Kazakhstanskaya Coal Mine, Kazakhstan, M1439, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- Sheet: Boundaries and methane sources --------------------------------
$data = $wb.Worksheets.Item("Boundaries and methane sources")

# Find the last used row in column B (the "id" column).
$lastRow = $data.Cells.Item($data.Rows.Count, 2).End(-4162).Row

# Collect the row numbers for the retired ids first, then delete from the
# bottom up so earlier row indices stay valid while deleting.
$idsToRemove = @("M1439.P2", "M1439.P5", "M1439.P7")
$rowsToDelete = @()

for ($r = 2; $r -le $lastRow; $r++) {
    $idVal = $data.Cells.Item($r, 2).Text
    if ($idsToRemove -contains $idVal) {
        $rowsToDelete += $r
    }
}

$rowsToDelete = $rowsToDelete | Sort-Object -Descending
foreach ($r in $rowsToDelete) {
    $data.Rows.Item($r).Delete()
}

# Refresh last row after deletions, then update the build_version column
# (S) for every remaining data row.
$lastRow = $data.Cells.Item($data.Rows.Count, 2).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $data.Cells.Item($r, 19).Value = $newVersion
}
